$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format before writing, so numeric-looking
# price strings (e.g. "1.003", "219.16") are preserved as text instead
# of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.191.76"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "1.658.57"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").Value = "219.16"
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("D6").Value = "0.5249"
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("D8").Value = "0.2674"
$ws.Range("E8").Value = "  +1.37%  "
$ws.Range("D9").Value = "0.06372"
$ws.Range("E9").Value = "  +0.77%  "
$ws.Range("D10").Value = "20.68"
$ws.Range("E10").Value = "  +0.12%  "
$ws.Range("D11").Value = "0.07696"
$ws.Range("E11").Value = "  -1.45%  "
$ws.Range("D12").Value = "4.658"
$ws.Range("E12").Value = "  +3.45%  "
$ws.Range("D13").Value = "1.649.89"
$ws.Range("E13").Value = "  -0.63%  "
$ws.Range("D14").Value = "1.887.03"
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("D15").Value = "0.5637"
$ws.Range("E15").Value = "  +1.35%  "
$ws.Range("D16").Value = "0.0₅8237"
$ws.Range("E16").Value = "  +2.64%  "
$ws.Range("D17").Value = "65.70"
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("D18").Value = "26.181.51"
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("E19").Value = "  -0.26%  "
$ws.Range("D20").Value = "4.679"
$ws.Range("E20").Value = "  +0.56%  "
$ws.Range("D21").Value = "10.52"
$ws.Range("E21").Value = "  +3.52%  "
$ws.Range("D22").Value = "192.84"
$ws.Range("E22").Value = "  -2.10%  "
$ws.Range("D23").Value = "5.989"
$ws.Range("E23").Value = "  +0.29%  "
$ws.Range("E24").Value = "  -0.29%  "
$ws.Range("D25").Value = "145.85"
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("D26").Value = "0.1202"
$ws.Range("E26").Value = "  -0.52%  "
$ws.Range("D27").Value = "7.316"
$ws.Range("E27").Value = "  +2.07%  "
$ws.Range("D28").Value = "16.03"
$ws.Range("E28").Value = "  -0.19%  "
$ws.Range("D29").Value = "1.524"
$ws.Range("E29").Value = "  +0.74%  "
$ws.Range("D30").Value = "0.05583"
$ws.Range("E30").Value = "  -3.39%  "
$ws.Range("D31").Value = "1.274"
$ws.Range("E31").Value = "  -0.35%  "
$ws.Range("D32").Value = "3.484"
$ws.Range("E32").Value = "  -0.26%  "
$ws.Range("D33").Value = "3.394"
$ws.Range("E33").Value = "  +1.32%  "
$ws.Range("D34").Value = "1.573"
$ws.Range("E34").Value = "  -0.81%  "
$ws.Range("D35").Value = "0.9550"
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").Value = "2.783"
$ws.Range("E36").Value = "  -0.97%  "
$ws.Range("E37").Value = "  -0.68%  "
$ws.Range("D38").Value = "0.5742"
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").Value = "0.01599"
$ws.Range("E39").Value = "  +0.27%  "
$ws.Range("D40").Value = "5.929"
$ws.Range("E40").Value = "  -0.73%  "
$ws.Range("E41").Value = "  -0.26%  "
$ws.Range("D42").Value = "1.034.66"
$ws.Range("E42").Value = "  -2.56%  "
$ws.Range("D43").Value = "0.8354"
$ws.Range("E43").Value = "  -2.11%  "
$ws.Range("D44").Value = "101.26"
$ws.Range("E44").Value = "  -1.81%  "
$ws.Range("D45").Value = "1.797.61"
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("D46").Value = "58.55"
$ws.Range("E46").Value = "  +0.44%  "
$ws.Range("D47").Value = "0.0₈106"
$ws.Range("E47").Value = "  +4.92%  "
$ws.Range("D48").Value = "1.000"
$ws.Range("E48").Value = "  -1.05%  "
$ws.Range("D49").Value = "8.077"
$ws.Range("E49").Value = "  +0.84%  "
$ws.Range("E50").Value = "  -1.41%  "
$ws.Range("D51").Value = "0.05239"
$ws.Range("E51").Value = "  +0.78%  "

# Restore the default cell style on column D (the text-format override
# was only needed to control how Excel parsed the values above).
$ws.Range("D2:D51").Style = "Normal"
